$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the descriptive/comment text cells in column C (and D2:D5) that were removed
$ws.Range("C1").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()

# Row 6: A6 changed from "dct:title" to "skos:prefLabel"
$ws.Range("A6").Value = "skos:prefLabel"
$ws.Range("C6").ClearContents()

$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
